$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 2516.6667
$ws.Range("I70").Value = 2350
$ws.Range("J70").Value = 2850
$ws.Range("K70").Value = 7050
$ws.Range("L70").Value = 8550
$ws.Range("M70").Value = -6780
$ws.Range("N70").Value = -9090
$ws.Range("H73").Value = 2516.6667
$ws.Range("I73").Value = 2350
$ws.Range("J73").Value = 2850
$ws.Range("K73").Value = 7050
$ws.Range("L73").Value = 8550
$ws.Range("M73").Value = -6114
$ws.Range("N73").Value = -10422
$ws.Range("H96").Value = 506.14285
$ws.Range("I96").Value = 175.16667
$ws.Range("J96").Value = 754.375
$ws.Range("K96").Value = 525.50001
$ws.Range("L96").Value = 2263.125
$ws.Range("M96").Value = 847.49999
$ws.Range("N96").Value = -5009.125
$ws.Range("H111").Value = 678.125
$ws.Range("I111").Value = 420.83334
$ws.Range("K111").Value = 1262.50002
$ws.Range("M111").Value = 1804.49998
$ws.Range("H112").Value = 1145.8334
$ws.Range("J112").Value = 1145.8334
$ws.Range("L112").Value = 3437.5002
$ws.Range("N112").Value = -5653.5002
$ws.Range("H118").Value = 325
$ws.Range("I118").Value = 325
$ws.Range("K118").Value = 975
$ws.Range("M118").Value = 682
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 793.55554
$ws.Range("I2").Value = 767.75
$ws.Range("K2").Value = 767.75
$ws.Range("M2").Value = -654.75
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("M31").ClearContents()
$ws.Range("H45").Value = 1291.6364
$ws.Range("I45").Value = 1226.25
$ws.Range("K45").Value = 1226.25
$ws.Range("M45").Value = -849.25
$ws.Range("H116").Value = 793.55554
$ws.Range("I116").Value = 767.75
$ws.Range("K116").Value = 767.75
$ws.Range("M116").Value = 1526.25
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").ClearContents()
$ws.Range("N133").Value = 0
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 793.55554
$ws.Range("I3").Value = 767.75
$ws.Range("K3").Value = 767.75
$ws.Range("M3").Value = -653.75
$ws.Range("H105").Value = 1499.75
$ws.Range("I105").Value = 1499.6666
$ws.Range("K105").Value = 1499.6666
$ws.Range("M105").Value = 247.3334
$ws.Range("H122").Value = 100000
$ws.Range("J122").Value = 100000
$ws.Range("L122").Value = 100000
$ws.Range("N122").Value = -109800
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 40117
$ws.Range("I31").Value = 28238.6
$ws.Range("K31").Value = 28238.6
$ws.Range("M31").Value = -27943.6
$ws.Range("H34").Value = 40117
$ws.Range("I34").Value = 28238.6
$ws.Range("K34").Value = 28238.6
$ws.Range("M34").Value = -28036.6
$ws.Range("H59").Value = 65000
$ws.Range("J59").Value = 65000
$ws.Range("L59").Value = 65000
$ws.Range("N59").Value = -67290
$ws.Range("H93").Value = 17500
$ws.Range("I93").Value = 17500
$ws.Range("K93").Value = 17500
$ws.Range("M93").Value = -15628
$ws.Range("H107").Value = 182.27777
$ws.Range("J107").Value = 109.4
$ws.Range("L107").Value = 109.4
$ws.Range("N107").Value = -3949.4
$ws.Range("H141").Value = 0
$ws.Range("I141").Value = 0
$ws.Range("K141").Value = 0
$ws.Range("M141").ClearContents()
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 1825
$ws.Range("I80").Value = 1499
$ws.Range("J80").Value = 2151
$ws.Range("K80").Value = 4497
$ws.Range("L80").Value = 6453
$ws.Range("M80").Value = -3561
$ws.Range("N80").Value = -8325
$ws.Range("H83").Value = 1825
$ws.Range("I83").Value = 1499
$ws.Range("J83").Value = 2151
$ws.Range("K83").Value = 13491
$ws.Range("L83").Value = 19359
$ws.Range("M83").Value = -8811
$ws.Range("N83").Value = -28719
$ws.Range("H97").Value = 194.63637
$ws.Range("I97").Value = 196.88889
$ws.Range("K97").Value = 590.6666700000001
$ws.Range("M97").Value = -94.66667000000007
$ws.Range("H113").Value = 280
$ws.Range("I113").Value = 280
$ws.Range("K113").Value = 840
$ws.Range("M113").Value = 1330
$ws.Range("H121").Value = 1477.6666
$ws.Range("I121").Value = 500
$ws.Range("J121").Value = 1966.5
$ws.Range("K121").Value = 1500
$ws.Range("L121").Value = 5899.5
$ws.Range("M121").Value = -190
$ws.Range("N121").Value = -8519.5
$ws.Range("H129").Value = 5333
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 5333
$ws.Range("K129").Value = 0
$ws.Range("L129").ClearContents()
$ws.Range("M129").Value = 15999
$ws.Range("N129").Value = -25999
$ws.Range("H131").Value = 2415.6667
$ws.Range("I131").Value = 1123.5
$ws.Range("J131").Value = 5000
$ws.Range("K131").Value = 3370.5
$ws.Range("L131").Value = 15000
$ws.Range("M131").Value = 1669.5
$ws.Range("N131").Value = -25080
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 1702.2
$ws.Range("I80").Value = 1552.5
$ws.Range("J80").Value = 1802
$ws.Range("K80").Value = 1552.5
$ws.Range("L80").Value = 1802
$ws.Range("M80").Value = -554.5
$ws.Range("N80").Value = -3798
$ws.Range("H83").Value = 1702.2
$ws.Range("I83").Value = 1552.5
$ws.Range("J83").Value = 1802
$ws.Range("K83").Value = 7762.5
$ws.Range("L83").Value = 9010
$ws.Range("M83").Value = -2770.5
$ws.Range("N83").Value = -18994
$ws.Range("H102").Value = 26184.625
$ws.Range("I102").Value = 27876.934
$ws.Range("K102").Value = 27876.934
$ws.Range("M102").Value = -26254.934
$ws.Range("H122").Value = 2709.9
$ws.Range("I122").Value = 2437.5
$ws.Range("J122").Value = 3799.5
$ws.Range("K122").Value = 7312.5
$ws.Range("L122").Value = 11398.5
$ws.Range("M122").Value = -4862.5
$ws.Range("N122").Value = -16298.5
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H60").Value = 45000
$ws.Range("J60").Value = 45000
$ws.Range("L60").Value = 45000
$ws.Range("N60").Value = -46018
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").ClearContents()
$ws.Range("N81").Value = 0
$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").ClearContents()
$ws.Range("N84").Value = 0
$ws.Range("H133").Value = 120000
$ws.Range("J133").Value = 120000
$ws.Range("L133").Value = 120000
$ws.Range("N133").Value = -125060
$ws.Range("H136").Value = 5467.6665
$ws.Range("I136").Value = 5221.2
$ws.Range("K136").Value = 15663.6
$ws.Range("M136").Value = -13113.6
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1200
$ws.Range("I96").Value = 1200
$ws.Range("K96").Value = 1200
$ws.Range("M96").Value = 173
$ws.Range("H122").Value = 1000.4286
$ws.Range("I122").Value = 900.8
$ws.Range("J122").Value = 1249.5
$ws.Range("K122").Value = 2702.4
$ws.Range("L122").Value = 3748.5
$ws.Range("M122").Value = -252.3999999999996
$ws.Range("N122").Value = -8648.5
$ws.Range("H126").Value = 2347.9
$ws.Range("I126").Value = 1998.625
$ws.Range("K126").Value = 5995.875
$ws.Range("M126").Value = -3525.875
$ws.Range("H136").Value = 6156.625
$ws.Range("I136").Value = 6156.625
$ws.Range("K136").Value = 18469.875
$ws.Range("M136").Value = -15919.875

Write-Output "Applied all changes"